# Friday Project 10 - update sentiment scores & categories
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scores = @(-0.09899030109339679,0.6475516089082036,0.04413380580529414,0.749322346792024,-0.6233702274888613,0.000878953486077666,0.2786825308053305,-0.3931367970387509,-0.05260177637738402,-0.406133221668862,0.00693790114760473,-0.4891867060569277,0.01484952414969645,0.01156747152330598,-0.5193943617369797,0.08611197800163442,-0.05612078416268695,0.3795056309923407,-0.4567682922686051,-0.7215454167070601,-0.704367386330657,0.05009646864856881,-0.08982086458710262,-0.2431844924257217,0.2330800645891585,0.08417874335569822,-0.2564881627220614,0.4084125241310997,-0.4061499691327267,0.4084053965437379,-0.5122847062040259,0.4702273907451405,0.03896718399218754,-0.6853712231951988,0.05639304090890654,-0.6957317975929354,-0.2208546591093898,-0.7045125830545405,-0.07021213215627867,-0.3188057290312843,-0.739159504672815,0.4455664616347744,-0.7262492821163337,-0.2090776473786516,0.03439923790126176,-0.386836377926267,-0.4213442211333231,0.04931711468551742,-0.3328429135022828,0.5876963779825463,-0.2462190137356823,-0.5984289697830769,-0.6917715232691061,0.7789985158936701,-0.353685057430468,0.5750940415225891,-0.7263934071703708,0.4526293770807616,0.07803519765759578,-0.2043378646821583,0.2734978812465426,-0.01835985961342608,-0.5698911478379212,-0.3983527429786172,-0.2676184988539407,-0.4889929247566807,0.08182282711265823,0.009915583424310356,-0.2883273487265714,0.2427487156843753,0.05553261524814987,-0.4123690346880222,0.669325149101067,-0.07570562440320004,-0.4235795808030751,-0.03730652901818653,0.006285009914272943,-0.7782500879418991,-0.5730518890819443,-0.4587741339993099)
$categories = @("Neutral","Positive","Neutral","Positive","Negative","Neutral","Positive","Negative","Neutral","Negative","Neutral","Negative","Neutral","Neutral","Negative","Neutral","Neutral","Positive","Negative","Negative","Negative","Neutral","Neutral","Negative","Positive","Neutral","Negative","Positive","Negative","Positive","Negative","Positive","Neutral","Negative","Neutral","Negative","Negative","Negative","Neutral","Negative","Negative","Positive","Negative","Negative","Neutral","Negative","Negative","Neutral","Negative","Positive","Negative","Negative","Negative","Positive","Negative","Positive","Negative","Positive","Neutral","Negative","Positive","Neutral","Negative","Negative","Negative","Negative","Neutral","Neutral","Negative","Positive","Neutral","Negative","Positive","Neutral","Negative","Neutral","Neutral","Negative","Negative","Negative")

for ($i = 0; $i -lt $scores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $scores[$i]
    $ws.Cells.Item($row, 4).Value = $categories[$i]
}
